$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / borders) from the "top of block" row (row 4)
# onto E7/F7, which are becoming the first row of a new time-entry block.
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new timesheet entry on row 7: date, start time, end time.
$ws.Range("B7").Value = 44758
$ws.Range("C7").Value = 0.90625
$ws.Range("D7").Value = 0.96875

# Give the new row 7 and the now-active row 8 their own "time diff" formula.
$ws.Range("E7").Formula = "=D7-C7"
$ws.Range("E8").Formula = "=D8-C8"

# Recompute the rolling-sum formulas for the smaller blocks.
$ws.Range("F4").Formula = "=SUM(E4:E6)"
$ws.Range("F7").Formula = "=SUM(E7:E10)"

# Update the remembered selection (where the user left the cursor).
$ws.Range("D7").Select() | Out-Null
